$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full refreshed data for the two CSV-backed columns (list size / running time
# in microseconds). Net effect vs. the previous snapshot: two new samples for
# list size 100 were captured (now rows 2-3) and the whole series shifted up
# by two rows, which also drops the oldest two samples that used to sit at
# the bottom of the 6553600 bucket (rows 172-173).
$listSize = @(100,100,100,100,100,100,100,100,100,100,200,200,200,200,200,200,200,200,200,200,400,400,400,400,400,400,400,400,400,400,800,800,800,800,800,800,800,800,800,800,1600,1600,1600,1600,1600,1600,1600,1600,1600,1600,3200,3200,3200,3200,3200,3200,3200,3200,3200,3200,6400,6400,6400,6400,6400,6400,6400,6400,6400,6400,12800,12800,12800,12800,12800,12800,12800,12800,12800,12800,25600,25600,25600,25600,25600,25600,25600,25600,25600,25600,51200,51200,51200,51200,51200,51200,51200,51200,51200,51200,102400,102400,102400,102400,102400,102400,102400,102400,102400,102400,204800,204800,204800,204800,204800,204800,204800,204800,204800,204800,409600,409600,409600,409600,409600,409600,409600,409600,409600,409600,819200,819200,819200,819200,819200,819200,819200,819200,819200,819200,1638400,1638400,1638400,1638400,1638400,1638400,1638400,1638400,1638400,1638400,3276800,3276800,3276800,3276800,3276800,3276800,3276800,3276800,3276800,3276800,6553600,6553600,6553600,6553600,6553600,6553600,6553600,6553600,6553600,6553600)
$runningSeconds = @(39370,36790,27132,42208,62877,90113,123619,159482,200890,261793,337458,381444,934874,541911,627374,701461,773248,876482,1021620,1043571,1181560,2219985,1540331,1514172,1524692,2260133,2497998,2043478,2154803,2385230,2744455,2815514,3502863,3427065,3026285,3165783,3528439,5276998,4405560,3801474,4323239,4488792,5781360,5735277,9488033,8512905,15400782,8650516,9763517,8290248,8760201,8976742,8193592,12286780,10249075,14857894,10651328,28877883,12951668,10028417,10256810,10396700,13093371,12448926,11447765,12141072,16142079,11762945,11866383,12426959,13766417,16248949,23669399,25100521,14207426,16402317,15686809,18971809,15633464,18319510,16646635,19046375,17609676,17794950,17999876,23330510,20581943,22230494,27728356,20790511,19947218,23017732,22730905,21466183,22847987,23881389,27858281,23377178,29251754,33492840,29643102,32106370,30381630,35825864,27916138,28583577,34257264,28895942,29290119,30415070,31596781,34440013,32024928,31954919,32866042,34746208,33801145,34064973,37027027,37217129,36992254,40837390,37926666,42009179,40435418,37816866,40805825,43220642,43911265,45365608,42512097,43677172,44633134,55006519,46161684,49289940,48520809,49948347,50184363,49774882,51995034,50635587,56273893,55102791,55929332,54401107,61554114,60512774,59928122,58157387,58336812,60181792,62078667,60283422,61677775,61205396,68895356,66282782,66385019,72784977,70089556,89107748,68820607,71510764,71585441,232099762,87873129,139086528,80789466,79802410)

for ($i = 0; $i -lt $listSize.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $listSize[$i]
    $ws.Cells.Item($i + 2, 2).Value = $runningSeconds[$i]
}

# The new series is two rows shorter than before (170 data rows instead of
# 172), so drop the now-stale trailing rows that still hold the old values.
$ws.Rows("172:173").Delete()

$ws.Activate()
$ws.Range("A1:B1").Select()

# The running-seconds query table now only covers columns C:FO (the first
# two data points moved under a separate capture), so repoint the named
# range accordingly.
$wb.Names("running_seconds_with_csv").RefersTo = "=Sheet1!`$C`$2:`$FO`$2"
